$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 18:30"

# Updated COVID-19 statistics per country (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 8675199
$ws.Range("C4").Value = 13548
$ws.Range("D4").Value = 5662998
$ws.Range("E4").Value = 2783624
$ws.Range("B5").Value = 7781746
$ws.Range("C5").Value = 22106
$ws.Range("D5").Value = 6965699
$ws.Range("E5").Value = 698573
$ws.Range("G5").Value = 138
$ws.Range("H5").Value = 117474
$ws.Range("E6").Value = 391341
$ws.Range("G6").Value = 34
$ws.Range("H6").Value = 155996
$ws.Range("B14").Value = 830998
$ws.Range("C14").Value = 20530
$ws.Range("G14").Value = 224
$ws.Range("H14").Value = 44571
$ws.Range("B20").Value = 414293
$ws.Range("C20").Value = 10419
$ws.Range("E20").Value = 94017
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = 10076
$ws.Range("B31").Value = 231146
$ws.Range("C31").Value = 8081
$ws.Range("D31").Value = 91589
$ws.Range("E31").Value = 137616
$ws.Range("G31").Value = 96
$ws.Range("H31").Value = 1941
$ws.Range("B32").Value = 228318
$ws.Range("C32").Value = 13632
$ws.Range("D32").Value = 105092
$ws.Range("E32").Value = 119054
$ws.Range("G32").Value = 153
$ws.Range("H32").Value = 4172
$ws.Range("B33").Value = 210879
$ws.Range("C33").Value = 1731
$ws.Range("D33").Value = 177297
$ws.Range("E33").Value = 23699
$ws.Range("G33").Value = 21
$ws.Range("H33").Value = 9883
$ws.Range("B49").Value = 103902
$ws.Range("C49").Value = 730
$ws.Range("D49").Value = 93341
$ws.Range("E49").Value = 6967
$ws.Range("G49").Value = 14
$ws.Range("H49").Value = 3594
$ws.Range("B50").Value = 103653
$ws.Range("C50").Value = 6634
$ws.Range("D50").Value = 55800
$ws.Range("E50").Value = 45790
$ws.Range("G50").Value = 11
$ws.Range("H50").Value = 2063
$ws.Range("B63").Value = 64724
$ws.Range("C63").Value = 285
$ws.Range("D63").Value = 61957
$ws.Range("E63").Value = 2225
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 542
$ws.Range("D65").Value = 57832
$ws.Range("E65").Value = 91
$ws.Range("B67").Value = 55630
$ws.Range("C67").Value = 273
$ws.Range("D67").Value = 38788
$ws.Range("E67").Value = 14945
$ws.Range("G67").Value = 9
$ws.Range("H67").Value = 1897
$ws.Range("B74").Value = 47843
$ws.Range("C74").Value = 631
$ws.Range("D74").Value = 33421
$ws.Range("E74").Value = 13538
$ws.Range("G74").Value = 14
$ws.Range("H74").Value = 884
$ws.Range("B75").Value = 47601
$ws.Range("C75").Value = 63
$ws.Range("D75").Value = 46824
$ws.Range("E75").Value = 463
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 314
$ws.Range("B76").Value = 47214
$ws.Range("C76").Value = 1322
$ws.Range("D76").Value = 5032
$ws.Range("E76").Value = 41398
$ws.Range("G76").Value = 44
$ws.Range("H76").Value = 784
$ws.Range("B87").Value = 29057
$ws.Range("C87").Value = 841
$ws.Range("E87").Value = 18509
$ws.Range("G87").Value = 10
$ws.Range("H87").Value = 559
$ws.Range("B99").Value = 16436
$ws.Range("C99").Value = 177
$ws.Range("D99").Value = 12378
$ws.Range("E99").Value = 3803
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 255
$ws.Range("B104").Value = 12851
$ws.Range("C104").Value = 518
$ws.Range("D104").Value = 9085
$ws.Range("E104").Value = 3625
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 141
$ws.Range("B105").Value = 12460
$ws.Range("D105").Value = 10609
$ws.Range("E105").Value = 1718
$ws.Range("H105").Value = 133
$ws.Range("B146").Value = 4038
$ws.Range("C146").Value = 227
$ws.Range("D146").Value = 2729
$ws.Range("E146").Value = 1240
$ws.Range("G146").Value = 6
$ws.Range("H146").Value = 69
$ws.Range("B147").Value = 3897
$ws.Range("D147").Value = 3166
$ws.Range("E147").Value = 629
$ws.Range("H147").Value = 102
$ws.Range("B148").Value = 3877
$ws.Range("D148").Value = 2853
$ws.Range("E148").Value = 907
$ws.Range("H148").Value = 117
$ws.Range("B154").Value = 2876
$ws.Range("C154").Value = 4
$ws.Range("E154").Value = 1530
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 56
$ws.Range("B177").Value = 553
$ws.Range("C177").Value = 2
$ws.Range("E177").Value = 55
$ws.Range("B204").Value = 37
$ws.Range("C204").Value = 4
$ws.Range("E204").Value = 8
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
